$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3.4165806442375592 * [Math]::Pow(10, -18)
$ws.Range("C3").Value = 1.2565254049682620 * [Math]::Pow(10, -17)
$ws.Range("C4").Value = 6.4230443534434707 * [Math]::Pow(10, -18)
$ws.Range("C5").Value = 5.6080164811450333 * [Math]::Pow(10, -19)
$ws.Range("C6").Value = 4.1057969038852444 * [Math]::Pow(10, -18)
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("C10").Value = 1.2121516977527739 * [Math]::Pow(10, -18)
$ws.Range("C11").Value = 0
$ws.Range("C12").Value = 3.1588837548770521 * [Math]::Pow(10, -18)
$ws.Range("C13").Value = 9.7554613671491473 * [Math]::Pow(10, -18)
$ws.Range("C14").Value = 5.7753613842574747 * [Math]::Pow(10, -1)
$ws.Range("C15").Value = 7.6359781954471548 * [Math]::Pow(10, -18)
$ws.Range("C16").Value = 4.6945720971355026 * [Math]::Pow(10, -18)
$ws.Range("C17").Value = 4.2246386157425259 * [Math]::Pow(10, -1)
